$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update symbols (column A, rows 2-5)
$ws.Range("A2").Value = "CYTO"
$ws.Range("A3").Value = "MLGO"
$ws.Range("A4").Value = "FLGC"
$ws.Range("A5").Value = "SMFL"

# Update dates (column B, rows 2-5) to 2024-03-25 (serial 45376)
$ws.Range("B2:B5").Value = Get-Date -Year 2024 -Month 3 -Day 25 -Hour 0 -Minute 0 -Second 0

# Remove row 6 entirely (used to hold HOLO / 45364)
$ws.Rows.Item(6).Delete()

# Update the selected cell shown in the saved view
$ws.Range("F18").Select()
